# Updates cryptos list values (price/volume) to match the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $cellRef, $newValue) {
    $range = $sheet.Range($cellRef)
    # Force text interpretation so numeric-looking strings (e.g. '315.00', '1.000')
    # keep their exact original formatting instead of being parsed into numbers.
    $range.NumberFormat = '@'
    $range.Value = $newValue
    # Restore the default (unstyled) cell style so no stray formatting is introduced.
    $range.Style = "Normal"
}

Set-TextValue $ws 'D2' '27.615.05'
Set-TextValue $ws 'E2' '  -2.42%  '
Set-TextValue $ws 'D3' '1.846.46'
Set-TextValue $ws 'E3' '  -1.36%  '
Set-TextValue $ws 'E4' '  -0.16%  '
Set-TextValue $ws 'D5' '315.00'
Set-TextValue $ws 'E6' '  -0.04%  '
Set-TextValue $ws 'D7' '0.4275'
Set-TextValue $ws 'E7' '  -2.90%  '
Set-TextValue $ws 'D8' '0.3654'
Set-TextValue $ws 'E8' '  -1.18%  '
Set-TextValue $ws 'D9' '45.75'
Set-TextValue $ws 'E9' '  +1.06%  '
Set-TextValue $ws 'D10' '0.07284'
Set-TextValue $ws 'E10' '  -3.16%  '
Set-TextValue $ws 'D11' '0.8966'
Set-TextValue $ws 'E11' '  -4.57%  '
Set-TextValue $ws 'D12' '20.70'
Set-TextValue $ws 'E12' '  -3.41%  '
Set-TextValue $ws 'D13' '1.854.88'
Set-TextValue $ws 'E13' '  -1.48%  '
Set-TextValue $ws 'D14' '5.395'
Set-TextValue $ws 'D15' '6.586'
Set-TextValue $ws 'E15' '  -2.23%  '
Set-TextValue $ws 'D16' '0.06863'
Set-TextValue $ws 'E16' '  +0.14%  '
Set-TextValue $ws 'D17' '1.001'
Set-TextValue $ws 'E17' '  -0.07%  '
Set-TextValue $ws 'D18' '78.62'
Set-TextValue $ws 'E18' '  -4.28%  '
Set-TextValue $ws 'D19' '0.000008888'
Set-TextValue $ws 'E19' '  -2.01%  '
Set-TextValue $ws 'D20' '1.000'
Set-TextValue $ws 'E20' '  -0.09%  '
Set-TextValue $ws 'D21' '15.56'
Set-TextValue $ws 'E21' '  -2.57%  '
Set-TextValue $ws 'D22' '27.621.20'
Set-TextValue $ws 'E22' '  -2.33%  '
Set-TextValue $ws 'D23' '4.993'
Set-TextValue $ws 'E23' '  -2.82%  '
Set-TextValue $ws 'E24' '  -2.88%  '
Set-TextValue $ws 'D25' '2.076.69'
Set-TextValue $ws 'E25' '  -0.81%  '
Set-TextValue $ws 'D26' '2.038'
Set-TextValue $ws 'E26' '  +0.52%  '
Set-TextValue $ws 'D27' '155.08'
Set-TextValue $ws 'E27' '  +0.13%  '
Set-TextValue $ws 'D28' '18.43'
Set-TextValue $ws 'E28' '  +0.19%  '
Set-TextValue $ws 'D29' '5.249'
Set-TextValue $ws 'E29' '  -1.84%  '
Set-TextValue $ws 'D30' '114.28'
Set-TextValue $ws 'E30' '  +0.47%  '
Set-TextValue $ws 'D31' '1.841'
Set-TextValue $ws 'E31' '  +6.22%  '
Set-TextValue $ws 'D32' '0.08901'
Set-TextValue $ws 'D33' '0.7855'
Set-TextValue $ws 'E33' '  -2.00%  '
Set-TextValue $ws 'E34' '  -5.59%  '
Set-TextValue $ws 'D35' '2.970'
Set-TextValue $ws 'E35' '  +2.14%  '
Set-TextValue $ws 'D36' '1.111'
Set-TextValue $ws 'E36' '  -5.64%  '
Set-TextValue $ws 'D37' '0.9999'
Set-TextValue $ws 'E37' '  -0.12%  '
Set-TextValue $ws 'D38' '0.05445'
Set-TextValue $ws 'E38' '  -0.03%  '
Set-TextValue $ws 'E39' '  -2.06%  '
Set-TextValue $ws 'D41' '2.776'
Set-TextValue $ws 'E41' '  -4.78%  '
Set-TextValue $ws 'D42' '0.5074'
Set-TextValue $ws 'E42' '  -3.59%  '
Set-TextValue $ws 'D43' '6.827'
Set-TextValue $ws 'E43' '  -4.16%  '
Set-TextValue $ws 'D44' '0.1651'
Set-TextValue $ws 'E44' '  -1.89%  '
Set-TextValue $ws 'D45' '8.265'
Set-TextValue $ws 'E45' '  -5.94%  '
Set-TextValue $ws 'D46' '0.06636'
Set-TextValue $ws 'E46' '  -1.82%  '
Set-TextValue $ws 'D47' '10.36'
Set-TextValue $ws 'E47' '  -1.50%  '
Set-TextValue $ws 'B48' 'Quant'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws 'D48' '105.86'
Set-TextValue $ws 'E48' '  -2.00%  '
Set-TextValue $ws 'B49' 'Decentraland'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws 'D49' '0.4718'
Set-TextValue $ws 'E49' '  -3.26%  '
Set-TextValue $ws 'D50' '1.000'
Set-TextValue $ws 'E50' '  -0.04%  '
Set-TextValue $ws 'D51' '1.637'
Set-TextValue $ws 'E51' '  -2.79%  '
